# Weekly roll-forward of the Zanahoria (carrot) price series.
# The sheet holds one row per market observation, ordered from most-recent
# (top, row 432 in the observed block) to oldest. A new week of data is
# published, so every existing observation shifts down by two rows
# (two new "Primera"/"Segunda" quality records are inserted at the top)
# and the two newest records are written at rows 432-433.
#
# Columns A-C are a constant template (market id / name / region) that is
# identical on every row; columns D-R carry the actual observation. The
# whole A-R span is shifted down by 2 rows (two brand-new rows appear at
# the bottom, 530-531, that also need A-C filled in), working from the
# bottom up so we never clobber a source row before it has been copied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 432
$lastDataRow = 529
$shift = 2

# Column letters A..R as 1-based column indices (A=1 .. R=18)
$firstCol = 1
$lastCol = 18

# 1) Snapshot the current D..R contents of rows 432..529 before writing
#    anything, so later writes never read already-overwritten data.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the snapshot back out shifted down by 2 rows: new row (r+2)
#    gets the old contents of row r. Go from the highest source row down
#    so the target rows are filled from the bottom of the range upward.
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $target = $r + $shift
    $rowVals = $snapshot[$r]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($target, $c).Value = $rowVals[$c]
    }
}

# 3) Populate the two new observations at the top of the block (rows
#    432-433) with this week's published data.

# Row 432: Primera quality
$ws.Range("D432").Value = 45211
$ws.Range("E432").Value = 16
$ws.Range("F432").Value = 100114013
$ws.Range("G432").Value = "Zanahoria"
$ws.Range("H432").Value = "Sin especificar"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 180
$ws.Range("K432").Value = 6500
$ws.Range("L432").Value = 7000
$ws.Range("M432").Value = 6778
$ws.Range("N432").Value = "$/saco 20 kilos"
$ws.Range("O432").Value = "Región de Ñuble"
$ws.Range("P432").Value = 339
$ws.Range("Q432").Value = 20
$ws.Range("R432").Value = "Hortaliza"

# Row 433: Segunda quality
$ws.Range("D433").Value = 45211
$ws.Range("E433").Value = 16
$ws.Range("F433").Value = 100114013
$ws.Range("G433").Value = "Zanahoria"
$ws.Range("H433").Value = "Sin especificar"
$ws.Range("I433").Value = "Segunda"
$ws.Range("J433").Value = 120
$ws.Range("K433").Value = 5500
$ws.Range("L433").Value = 6000
$ws.Range("M433").Value = 5750
$ws.Range("N433").Value = "$/saco 20 kilos"
$ws.Range("O433").Value = "Región de Ñuble"
$ws.Range("P433").Value = 288
$ws.Range("Q433").Value = 20
$ws.Range("R433").Value = "Hortaliza"
